# WeaponData.xlsx - add FireBall weapon (data row + dedicated sheet)

$wb = $excel.ActiveWorkbook
$wsWeaponData = $wb.Worksheets.Item("WeaponData")
$wsMagicBall  = $wb.Worksheets.Item("MagicBall")

# --- 1. WeaponData sheet: MagicBall's NowWeaponLevel drops from 1 to 0 ---
$wsWeaponData.Range("D4").Value = 0

# --- 2. WeaponData sheet: fill in the new FireBall row (row 5) ---
$wsWeaponData.Range("A5").Value = "FireBall"
$wsWeaponData.Range("B5").Value = "Assets/Prefabs/Weapons/FireBall.prefab"
$wsWeaponData.Range("C5").Value = "FireBall"
$wsWeaponData.Range("D5").Value = 0
$wsWeaponData.Range("E5").Value = "Assets/ArtResources/Weapons/RotateFire/Effect3/1.png"
$wsWeaponData.Range("F5").Value = 1
$wsWeaponData.Range("G5").Value = "周囲を回転しながら攻撃します。"
$wsWeaponData.Range("H5").Value = "FireBallController"

# C5 / F5 are brand-new cells (previously blank / absent) - match the
# formatting already used by the rest of the row.
$wsWeaponData.Range("A5").Copy()
$wsWeaponData.Range("C5").PasteSpecial(-4122)
$wsWeaponData.Range("D5").Copy()
$wsWeaponData.Range("F5").PasteSpecial(-4122)

# --- 3. Add a new "FireBall" sheet, positioned after "MagicBall" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFireBall = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsFireBall.Name = "FireBall"

# Mirror MagicBall's layout/formatting (same columns, same stat table shape).
$wsMagicBall.Range("A1:F7").Copy()
$wsFireBall.Range("A1:F7").PasteSpecial(-4122)
$wsFireBall.Columns.Item(5).ColumnWidth = $wsMagicBall.Columns.Item(5).ColumnWidth

# Headers (english / localized / type row) are identical to MagicBall's.
$wsFireBall.Range("A1").Value = "Level"
$wsFireBall.Range("B1").Value = "Hurt"
$wsFireBall.Range("C1").Value = "Number"
$wsFireBall.Range("D1").Value = "Speed"
$wsFireBall.Range("E1").Value = "Cooldown"
$wsFireBall.Range("F1").Value = "Duration"

$wsFireBall.Range("A2").Value = "階段"
$wsFireBall.Range("B2").Value = "傷害"
$wsFireBall.Range("C2").Value = "數量"
$wsFireBall.Range("D2").Value = "速度"
$wsFireBall.Range("E2").Value = "冷卻"
$wsFireBall.Range("F2").Value = "持續時間"

$wsFireBall.Range("A3").Value = "int"
$wsFireBall.Range("B3").Value = "float"
$wsFireBall.Range("C3").Value = "int"
$wsFireBall.Range("D3").Value = "float"
$wsFireBall.Range("E3").Value = "float"
$wsFireBall.Range("F3").Value = "float"

# FireBall's own level-up stat table.
$wsFireBall.Range("A4").Value = 1
$wsFireBall.Range("B4").Value = 10
$wsFireBall.Range("C4").Value = 1
$wsFireBall.Range("D4").Value = 1
$wsFireBall.Range("E4").Value = 10
$wsFireBall.Range("F4").Value = 5

$wsFireBall.Range("A5").Value = 2
$wsFireBall.Range("B5").Value = 10
$wsFireBall.Range("C5").Value = 3
$wsFireBall.Range("D5").Value = 2
$wsFireBall.Range("E5").Value = 11.5
$wsFireBall.Range("F5").Value = 7.5

$wsFireBall.Range("A6").Value = 3
$wsFireBall.Range("B6").Value = 40
$wsFireBall.Range("C6").Value = 3
$wsFireBall.Range("D6").Value = 4
$wsFireBall.Range("E6").Value = 13
$wsFireBall.Range("F6").Value = 10

$wsFireBall.Range("A7").Value = 4
$wsFireBall.Range("B7").Value = 60
$wsFireBall.Range("C7").Value = 5
$wsFireBall.Range("D7").Value = 8
$wsFireBall.Range("E7").Value = 14.5
$wsFireBall.Range("F7").Value = 12.5
